$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Mohsen Saeed Abadi's degree / scores to row 13
$ws.Range("F13").Value = 93
$ws.Range("G13").Value = 100
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("K13").Value = "68 70per"
$ws.Range("J13").Value = "76.5 90per"

# Update the active view/selection to reflect the latest edit position
$ws.Activate()
$ws.Range("J15").Select()
